$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "AddCustomerTest"

$ws.Range("A2").Value = "Alice"
$ws.Range("B2").Value = "Jonhson"

$ws.Range("A1").Value = "fistName"
$ws.Range("B1").Value = "lastName"
$ws.Range("C1").Value = "postCode"
$ws.Range("D1").Value = "alertText"

$ws.Range("C2").Value = 11230
$ws.Range("D2").Value = "Customer added successfully"

$ws.Range("D2").Select()
